$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32, shifting existing rows 32:45 down to 33:46
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new data point
$ws.Range("A32").Value = 10
$ws.Range("B32").Value = "Vega Modelo de Temuco"
$ws.Range("C32").Value = "La Araucanía"
$ws.Range("D32").Value = Get-Date -Year 2022 -Month 9 -Day 21 -Hour 0 -Minute 0 -Second 0
$ws.Range("E32").Value = 9
$ws.Range("F32").Value = 100112042
$ws.Range("G32").Value = "Locoto"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 30
$ws.Range("K32").Value = 2700
$ws.Range("L32").Value = 2700
$ws.Range("M32").Value = 2700
$ws.Range("N32").Value = "$/kilo"
$ws.Range("O32").Value = "Región de Arica y Parinacota"
$ws.Range("P32").Value = 2700
$ws.Range("Q32").Value = 1
$ws.Range("R32").Value = "Hortaliza"
